$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9000.333000000001
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 11000.5
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 11000.5
$ws.Range("M40").Value = -4825
$ws.Range("N40").Value = -11350.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4525.091
$ws.Range("I86").Value = 2977.5454
$ws.Range("J86").Value = 6072.636
$ws.Range("K86").Value = 2977.5454
$ws.Range("L86").Value = 6072.636
$ws.Range("M86").Value = -1854.5454
$ws.Range("N86").Value = -8318.636

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4525.091
$ws.Range("I89").Value = 2977.5454
$ws.Range("J89").Value = 6072.636
$ws.Range("K89").Value = 14887.727
$ws.Range("L89").Value = 30363.18
$ws.Range("M89").Value = -9271.726999999999
$ws.Range("N89").Value = -41595.18

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2749.1177
$ws.Range("I106").Value = 1561.3334
$ws.Range("K106").Value = 1561.3334
$ws.Range("M106").Value = -930.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 12823146
$ws.Range("I137").Value = 38463110
$ws.Range("K137").Value = 115389330
$ws.Range("M137").Value = -115386780

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2912.2021
$ws.Range("I138").Value = 2087.353
$ws.Range("J138").Value = 3094.3118
$ws.Range("K138").Value = 6262.059
$ws.Range("L138").Value = 9282.9354
$ws.Range("M138").Value = -1122.059
$ws.Range("N138").Value = -19562.9354

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4981.9287
$ws.Range("I2").Value = 971.7143
$ws.Range("K2").Value = 971.7143
$ws.Range("M2").Value = -858.7143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 22500
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 8000
$ws.Range("I55").Value = 8000
$ws.Range("K55").Value = 8000
$ws.Range("M55").Value = -7685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 15876280
$ws.Range("I74").Value = 30306030
$ws.Range("K74").Value = 30306030
$ws.Range("M74").Value = -30305156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 15876280
$ws.Range("I77").Value = 30306030
$ws.Range("K77").Value = 151530150
$ws.Range("M77").Value = -151525782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1134.2354
$ws.Range("I97").Value = 1121.8125
$ws.Range("K97").Value = 1121.8125
$ws.Range("M97").Value = -625.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 39999.5
$ws.Range("J101").Value = 39999.5
$ws.Range("L101").Value = 39999.5
$ws.Range("N101").Value = -46489.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4981.9287
$ws.Range("I116").Value = 971.7143
$ws.Range("K116").Value = 971.7143
$ws.Range("M116").Value = 1322.2857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3618
$ws.Range("I122").Value = 3020.08
$ws.Range("K122").Value = 9060.24
$ws.Range("M122").Value = -6610.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 43746.57
$ws.Range("J135").Value = 43746.57
$ws.Range("L135").Value = 43746.57
$ws.Range("N135").Value = -53886.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4981.9287
$ws.Range("I3").Value = 971.7143
$ws.Range("K3").Value = 971.7143
$ws.Range("M3").Value = -857.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3249.75
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3999.6667
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 3999.6667
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -5995.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 3249.75
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3999.6667
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 19998.3335
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -29982.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1885.875
$ws.Range("I134").Value = 1193.0454
$ws.Range("K134").Value = 3579.1362
$ws.Range("M134").Value = -1044.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1395.4546
$ws.Range("I16").Value = 759.13336
$ws.Range("K16").Value = 759.13336
$ws.Range("M16").Value = -472.13336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8853.286
$ws.Range("I22").Value = 6991.3335
$ws.Range("K22").Value = 6991.3335
$ws.Range("M22").Value = -6641.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2271.9
$ws.Range("I105").Value = 873.0769
$ws.Range("K105").Value = 873.0769
$ws.Range("M105").Value = 873.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1395.4546
$ws.Range("I113").Value = 759.13336
$ws.Range("K113").Value = 759.13336
$ws.Range("M113").Value = 1410.86664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35563876
$ws.Range("I4").Value = 45833520
$ws.Range("K4").Value = 137500560
$ws.Range("M4").Value = -137500448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 385000.7
$ws.Range("J37").Value = 385000.7
$ws.Range("L37").Value = 1155002.1
$ws.Range("N37").Value = -1155226.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 121.125
$ws.Range("I38").Value = 163.8
$ws.Range("K38").Value = 491.4
$ws.Range("M38").Value = -144.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2162.9167
$ws.Range("I46").Value = 406.875
$ws.Range("J46").Value = 3040.9375
$ws.Range("K46").Value = 1220.625
$ws.Range("L46").Value = 9122.8125
$ws.Range("M46").Value = -1129.625
$ws.Range("N46").Value = -9304.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 10014
$ws.Range("J69").Value = 10014
$ws.Range("L69").Value = 30042
$ws.Range("N69").Value = -31664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 10014
$ws.Range("J72").Value = 10014
$ws.Range("L72").Value = 90126
$ws.Range("N72").Value = -98238

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 166672800
$ws.Range("J75").Value = 8753.5
$ws.Range("L75").Value = 26260.5
$ws.Range("N75").Value = -28256.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 166672800
$ws.Range("J78").Value = 8753.5
$ws.Range("L78").Value = 78781.5
$ws.Range("N78").Value = -88765.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1150.2222
$ws.Range("I121").Value = 1197.4
$ws.Range("J121").Value = 1091.25
$ws.Range("K121").Value = 3592.2
$ws.Range("L121").Value = 3273.75
$ws.Range("M121").Value = -2282.2
$ws.Range("N121").Value = -5893.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5529.294
$ws.Range("J122").Value = 7681.125
$ws.Range("L122").Value = 23043.375
$ws.Range("N122").Value = -27943.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3725.1177
$ws.Range("J132").Value = 10895.5
$ws.Range("L132").Value = 32686.5
$ws.Range("N132").Value = -37746.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8161.421
$ws.Range("I46").Value = 4070.625
$ws.Range("J46").Value = 11136.546
$ws.Range("K46").Value = 4070.625
$ws.Range("L46").Value = 11136.546
$ws.Range("M46").Value = -3882.625
$ws.Range("N46").Value = -11512.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1517288.8
$ws.Range("I55").Value = 3334772
$ws.Range("J55").Value = 2719.4443
$ws.Range("K55").Value = 3334772
$ws.Range("L55").Value = 2719.4443
$ws.Range("M55").Value = -3334599
$ws.Range("N55").Value = -3065.4443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 50000
$ws.Range("K63").Value = 50000
$ws.Range("M63").Value = -49251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 50000
$ws.Range("K66").Value = 150000
$ws.Range("M66").Value = -146256

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 450178.44
$ws.Range("J122").Value = 7750.5
$ws.Range("L122").Value = 23251.5
$ws.Range("N122").Value = -28151.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 67565.5
$ws.Range("J109").Value = 67565.5
$ws.Range("L109").Value = 67565.5
$ws.Range("N109").Value = -70339.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3359.2222
$ws.Range("I122").Value = 1980.0769
$ws.Range("J122").Value = 6945
$ws.Range("K122").Value = 5940.2307
$ws.Range("L122").Value = 20835
$ws.Range("M122").Value = -3490.2307
$ws.Range("N122").Value = -25735
